$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record (Plátano, Vega Monumental Concepción) was
# reported; insert it as a new row ahead of the existing history so the
# sheet stays ordered, pushing the prior rows (and the former last row)
# down by one.
$ws.Rows.Item(334).Insert()

$ws.Range("A334").Value = 11
$ws.Range("B334").Value = "Vega Monumental Concepción"
$ws.Range("C334").Value = "Bíobío"
$ws.Range("D334").Value = 44588
$ws.Range("E334").Value = 8
$ws.Range("F334").Value = "Fruta"
$ws.Range("G334").Value = 100108
$ws.Range("H334").Value = "Tropicales y subtropicales"
$ws.Range("I334").Value = 100108006
$ws.Range("J334").Value = "Plátano"
$ws.Range("K334").Value = "Sin especificar"
$ws.Range("L334").Value = "Pintón"
$ws.Range("M334").Value = 220
$ws.Range("N334").Value = 11000
$ws.Range("O334").Value = 12000
$ws.Range("P334").Value = 11455
$ws.Range("Q334").Value = "$/caja 20 kilos"
$ws.Range("R334").Value = "Ecuador"
$ws.Range("S334").Value = 573
$ws.Range("T334").Value = 20
